$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the combined "Name Surname" column into three columns:
# insert two new blank columns before column B (old B/C shift to D/E, etc.)
$ws.Columns("B:C").Insert()

# New header for the first inserted column: title/prefix (คำนำหน้า)
$ws.Range("B1").Value = "คำนำหน้า"

# New header for the second inserted column: first name (ชื่อ)
$ws.Range("C1").Value = "ชื่อ"

# The old "ชื่อ นามสกุล" header (now shifted to D1) becomes just "นามสกุล" (surname)
$ws.Range("D1").Value = "นามสกุล"

# Update the active selection to reflect the new header row range
$ws.Range("A2:F2").Select() | Out-Null
